# Add a "Save" column (column H) to the s_vals sheet.
#
# Column H gets:
#   - H1: header label "Save", formatted like the other header cells
#         (bold, centered, bordered - i.e. same style as G1)
#   - H2:H4: numeric 0 values for each data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the neighboring "sum" header (G1) onto H1
# so the new header cell reuses the existing header style instead of Excel
# minting a brand new (duplicate) cell style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" column values for the existing data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
